$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.973.41'
$ws.Range("E2").Value = '  +1.99%  '
$ws.Range("D3").Value = '1.815.68'
$ws.Range("E3").Value = '  +2.51%  '
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = "'312.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.35%  '
$ws.Range("D7").Value = "'0.4297"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").Value = "'0.3666"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = "'0.07262"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("D10").Value = '2.160.86'
$ws.Range("E10").Value = '  +21.67%  '
$ws.Range("D11").Value = "'0.8639"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'21.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.77%  '
$ws.Range("D13").Value = "'5.410"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'6.599"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.51%  '
$ws.Range("D15").Value = "'0.06951"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.44%  '
$ws.Range("D16").Value = "'81.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.15%  '
$ws.Range("E17").Value = '  +0.72%  '
$ws.Range("D18").Value = "'0.000008889"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.67%  '
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").Value = "'15.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.79%  '
$ws.Range("D21").Value = '27.033.30'
$ws.Range("E21").Value = '  +2.24%  '
$ws.Range("D22").Value = "'5.177"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.46%  '
$ws.Range("D23").Value = '2.424.47'
$ws.Range("E23").Value = '  +21.53%  '
$ws.Range("D24").Value = "'11.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.30%  '
$ws.Range("D25").Value = "'153.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("D26").Value = "'1.885"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.93%  '
$ws.Range("D27").Value = "'18.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("D28").Value = "'5.224"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.40%  '
$ws.Range("D29").Value = "'1.897"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.87%  '
$ws.Range("D30").Value = "'114.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = "'0.08933"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.42%  '
$ws.Range("E32").Value = '  +6.10%  '
$ws.Range("D33").Value = "'0.7460"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.86%  '
$ws.Range("D34").Value = "'4.415"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'2.805"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.22%  '
$ws.Range("E36").Value = '  +0.42%  '
$ws.Range("E37").Value = '  +4.89%  '
$ws.Range("D38").Value = "'0.05207"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.76%  '
$ws.Range("D39").Value = "'0.01922"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("D40").Value = "'0.5099"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.39%  '
$ws.Range("E41").Value = '  +3.00%  '
$ws.Range("D42").Value = "'2.739"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.01%  '
$ws.Range("D43").Value = "'6.444"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.49%  '
$ws.Range("D44").Value = "'8.303"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.46%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = "'106.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'10.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.96%  '
$ws.Range("D47").Value = "'1.005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").Value = "'0.4581"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.20%  '
$ws.Range("D49").Value = "'1.643"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.52%  '
$ws.Range("D50").Value = "'0.06210"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("D51").Value = "'1.838"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.41%  '
